$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark sitting alone in the very first
#    (empty) paragraph of the document -> becomes a plain empty paragraph.
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$firstRange = $firstPara.Range
$xmlEmptyPara = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$firstRange.InsertXML($xmlEmptyPara)

# ---------------------------------------------------------------------------
# 2) Re-add the "_GoBack" bookmark at the end of the paragraph that talks
#    about the default 1280x720 window size (right before the paragraph
#    mark, after the last run).
# ---------------------------------------------------------------------------
$find = $d.Content
$found = $find.Find.Execute("Par défaut, les dimensions sont fixées", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 1280x720 paragraph"
}
$find.Expand(4) | Out-Null

$xmlParaWithBookmark = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Titre3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:firstLine="432"/></w:pPr><w:r><w:t xml:space="preserve">La fenêtre de configuration </w:t></w:r><w:r><w:t xml:space="preserve">(voir en figure 1) </w:t></w:r><w:r><w:t>va vous permettre de configurer comme vous le désirez votre simulateur. Vous pourrez décider de munir votre distributeur des composants que vous souhaitez. En supplément, vous avez la possibilité de choisir la taille de la fenêtre de simulation comme bon vous semble. Grâce à cela, vous pourrez tester la taille de fenêtre qui convient le mieux pour une telle application.</w:t></w:r><w:r><w:t xml:space="preserve"> Par défaut, les dimensions sont fixées à 1280x720. Vous ne devez donc pas forcément taper les dimensions à chaque fois que vous utilisez l’application.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$find.InsertXML($xmlParaWithBookmark)

# ---------------------------------------------------------------------------
# 3) Replace the empty paragraph that follows the last table (right before
#    "Affichage des horaires de train") with a small separator paragraph and
#    a new "N.B. : ..." note paragraph about the national registry number.
# ---------------------------------------------------------------------------
$lastTable = $d.Tables.Item($d.Tables.Count)
$afterTable = $lastTable.Range.Next(4, 1)  # wdParagraph = 4, count = 1 -> paragraph right after the table

$xmlNewParas = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="8"/><w:szCs w:val="8"/><w:u w:val="single"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>N.B.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>: Le registre national est accepté si et seulement s</w:t></w:r><w:r><w:t>’i</w:t></w:r><w:r><w:t>l comporte exactement 11 chiffres.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$afterTable.InsertXML($xmlNewParas)

Write-Host "done"
